# Auto-generated script to update cryptos worksheet values
# per commit "Updated cryptos list on Wed Aug 23 16:33:17 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.324.42'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.674.40'
$ws.Range("E3").Value = '  +2.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  -0.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.20'
$ws.Range("E5").Value = '  +5.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5296'
$ws.Range("E6").Value = '  +2.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9986'
$ws.Range("E7").Value = '  -0.65%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2653'
$ws.Range("E8").Value = '  +3.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06441'
$ws.Range("E9").Value = '  +3.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.21'
$ws.Range("E10").Value = '  +2.83%  '
$ws.Range("E11").Value = '  +3.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.682.19'
$ws.Range("E12").Value = '  +2.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.489'
$ws.Range("E13").Value = '  +2.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.900.33'
$ws.Range("E14").Value = '  +1.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5577'
$ws.Range("E15").Value = '  +4.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8396'
$ws.Range("E16").Value = '  +6.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.81'
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.356.91'
$ws.Range("E18").Value = '  +1.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9991'
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.809'
$ws.Range("E20").Value = '  +3.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '194.42'
$ws.Range("E21").Value = '  +4.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.34'
$ws.Range("E22").Value = '  +4.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.395'
$ws.Range("E23").Value = '  +4.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.98'
$ws.Range("E25").Value = '  -3.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1265'
$ws.Range("E26").Value = '  +4.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.468'
$ws.Range("E27").Value = '  +2.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.21'
$ws.Range("E28").Value = '  +4.29%  '
$ws.Range("E29").Value = '  +3.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06166'
$ws.Range("E30").Value = '  +3.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.274'
$ws.Range("E31").Value = '  +2.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.543'
$ws.Range("E32").Value = '  +3.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.454'
$ws.Range("E33").Value = '  +2.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.694'
$ws.Range("E34").Value = '  +4.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.010'
$ws.Range("E35").Value = '  +4.43%  '
$ws.Range("B36").Value = 'MXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.783'
$ws.Range("E36").Value = '  +2.24%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.399'
$ws.Range("E37").Value = '  +0.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5734'
$ws.Range("E38").Value = '  -1.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01624'
$ws.Range("E39").Value = '  +2.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.921'
$ws.Range("E40").Value = '  -0.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8648'
$ws.Range("E41").Value = '  +2.74%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.045.29'
$ws.Range("E42").Value = '  -3.40%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9982'
$ws.Range("E43").Value = '  -0.56%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.84'
$ws.Range("E44").Value = '  -0.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.820.31'
$ws.Range("E45").Value = '  +1.23%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.82'
$ws.Range("E46").Value = '  +4.43%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₈106'
$ws.Range("E47").Value = '  -3.64%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.177'
$ws.Range("E48").Value = '  +2.84%  '
$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.002'
$ws.Range("E49").Value = '  -0.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05191'
$ws.Range("E50").Value = '  -0.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.038'
$ws.Range("E51").Value = '  +3.95%  '
